# Updated symbol list - applies numeric/percent/rank updates to cryptos sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "306.75"
Set-TextValue $ws.Range("E2") "5.94%"
Set-TextValue $ws.Range("G2") "5"

Set-TextValue $ws.Range("D3") "32.19"
Set-TextValue $ws.Range("E3") "8.98%"
Set-TextValue $ws.Range("G3") "5"

Set-TextValue $ws.Range("D4") "5.335"
Set-TextValue $ws.Range("E4") "4.20%"
Set-TextValue $ws.Range("G4") "5"

Set-TextValue $ws.Range("D5") "0.07455"
Set-TextValue $ws.Range("G5") "5"

Set-TextValue $ws.Range("D6") "7.743"
Set-TextValue $ws.Range("E6") "5.21%"
Set-TextValue $ws.Range("G6") "5"

Set-TextValue $ws.Range("D7") "3.697"
Set-TextValue $ws.Range("E7") "8.53%"
Set-TextValue $ws.Range("G7") "5"

Set-TextValue $ws.Range("D8") "1.588"
Set-TextValue $ws.Range("E8") "17.30%"
Set-TextValue $ws.Range("G8") "5"

Set-TextValue $ws.Range("D9") "0.9203"
Set-TextValue $ws.Range("E9") "-0.02%"
Set-TextValue $ws.Range("G9") "5"

Set-TextValue $ws.Range("D10") "0.01644"
Set-TextValue $ws.Range("E10") "2,441.25%"
Set-TextValue $ws.Range("G10") "5"

Set-TextValue $ws.Range("D11") "0.1675"
Set-TextValue $ws.Range("E11") "5.47%"
Set-TextValue $ws.Range("G11") "5"

Set-TextValue $ws.Range("D12") "0.07666"
Set-TextValue $ws.Range("E12") "14.95%"
Set-TextValue $ws.Range("G12") "5"

Set-TextValue $ws.Range("D13") "0.07989"
Set-TextValue $ws.Range("E13") "3.11%"
Set-TextValue $ws.Range("G13") "5"

Set-TextValue $ws.Range("D14") "0.03078"
Set-TextValue $ws.Range("E14") "4.74%"
Set-TextValue $ws.Range("G14") "5"

Set-TextValue $ws.Range("D15") "0.09863"
Set-TextValue $ws.Range("E15") "9.75%"
Set-TextValue $ws.Range("G15") "5"

Set-TextValue $ws.Range("D16") "0.001531"
Set-TextValue $ws.Range("E16") "-3.36%"
Set-TextValue $ws.Range("G16") "5"

Set-TextValue $ws.Range("D17") "0.04563"
Set-TextValue $ws.Range("E17") "0.94%"
Set-TextValue $ws.Range("G17") "5"

Set-TextValue $ws.Range("D18") "0.006459"
Set-TextValue $ws.Range("E18") "3.14%"
Set-TextValue $ws.Range("G18") "5"

Set-TextValue $ws.Range("D19") "3.465"
Set-TextValue $ws.Range("E19") "0.41%"
Set-TextValue $ws.Range("G19") "5"

Set-TextValue $ws.Range("D20") "2.243"
Set-TextValue $ws.Range("E20") "1.08%"
Set-TextValue $ws.Range("G20") "5"

Set-TextValue $ws.Range("E21") "2.16%"
Set-TextValue $ws.Range("G21") "5"

Set-TextValue $ws.Range("D22") "0.1317"
Set-TextValue $ws.Range("E22") "0.59%"
Set-TextValue $ws.Range("G22") "5"

Set-TextValue $ws.Range("D23") "4.219"
Set-TextValue $ws.Range("E23") "3.88%"
Set-TextValue $ws.Range("G23") "5"

Set-TextValue $ws.Range("E24") "3.97%"
Set-TextValue $ws.Range("G24") "5"

Set-TextValue $ws.Range("E25") "2.13%"
Set-TextValue $ws.Range("G25") "5"

Set-TextValue $ws.Range("D26") "0.004530"
Set-TextValue $ws.Range("E26") "9.67%"
Set-TextValue $ws.Range("G26") "5"

Set-TextValue $ws.Range("D27") "0.0001168"
Set-TextValue $ws.Range("E27") "-6.44%"
Set-TextValue $ws.Range("G27") "5"

Set-TextValue $ws.Range("D28") "0.0001739"
Set-TextValue $ws.Range("E28") "7.56%"
Set-TextValue $ws.Range("G28") "5"

Set-TextValue $ws.Range("G29") "5"

Set-TextValue $ws.Range("G30") "5"

Set-TextValue $ws.Range("G31") "5"

Set-TextValue $ws.Range("G32") "5"

Set-TextValue $ws.Range("G33") "5"

Set-TextValue $ws.Range("G34") "5"

Set-TextValue $ws.Range("G35") "5"

Set-TextValue $ws.Range("G36") "5"

Set-TextValue $ws.Range("G37") "5"

Set-TextValue $ws.Range("G38") "5"

Set-TextValue $ws.Range("G39") "5"

Set-TextValue $ws.Range("D40") "0.04521"
Set-TextValue $ws.Range("E40") "6.98%"
Set-TextValue $ws.Range("G40") "5"

Set-TextValue $ws.Range("D41") "0.007426"
Set-TextValue $ws.Range("E41") "10.30%"
Set-TextValue $ws.Range("G41") "5"

Set-TextValue $ws.Range("D42") "0.1366"
Set-TextValue $ws.Range("E42") "10.05%"
Set-TextValue $ws.Range("G42") "5"

Set-TextValue $ws.Range("D43") "0.002257"
Set-TextValue $ws.Range("E43") "14.08%"
Set-TextValue $ws.Range("G43") "5"

Set-TextValue $ws.Range("D44") "0.01369"
Set-TextValue $ws.Range("E44") "6.45%"
Set-TextValue $ws.Range("G44") "5"

Set-TextValue $ws.Range("D45") "0.00006127"
Set-TextValue $ws.Range("E45") "9.86%"
Set-TextValue $ws.Range("G45") "5"

Set-TextValue $ws.Range("D46") "1.893"
Set-TextValue $ws.Range("E46") "-3.93%"
Set-TextValue $ws.Range("G46") "5"

Set-TextValue $ws.Range("D47") "0.01298"
Set-TextValue $ws.Range("E47") "-0.57%"
Set-TextValue $ws.Range("G47") "5"

Set-TextValue $ws.Range("G48") "5"

Set-TextValue $ws.Range("G49") "5"

Set-TextValue $ws.Range("G50") "5"

Set-TextValue $ws.Range("G51") "5"

